$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE  003 -> 001  (must stay text, not become the number 1)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# M2: NOTICE_DATE (stored as text)
$ws.Range("M2").Value = "2020-12-22 00:00:00"

# N2: REPORT_DATE (stored as text)
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# O2 / P2: NETCASH_OPERATE & its ratio
$ws.Range("O2").Value = 45685238.46
$ws.Range("P2").Value = 65.51186715670001

# Q2 / R2: SALES_SERVICES & ratio
$ws.Range("Q2").Value = 133435570.18
$ws.Range("R2").Value = 191.3443738564

# S2 / T2: PAY_STAFF_CASH & ratio
$ws.Range("S2").Value = 10855145.24
$ws.Range("T2").Value = 15.5660965533

# U2 / V2: NETCASH_INVEST & ratio
$ws.Range("U2").Value = -1613589.13
$ws.Range("V2").Value = -2.3138598001

# Y2 / Z2: CONSTRUCT_LONG_ASSET & ratio
$ws.Range("Y2").Value = 3404043.09
$ws.Range("Z2").Value = 4.8813408056

# AA2 / AB2: NETCASH_FINANCE & ratio (previously blank, now populated)
$ws.Range("AA2").Value = 25664800
$ws.Range("AB2").Value = 36.8028935577

# AC2 / AD2: CCE_ADD & ratio
$ws.Range("AC2").Value = 69735821.06999999
$ws.Range("AD2").Value = 23.157611151
